$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (this also updates the Print_Area defined name and the
# workbook's internal references automatically, since they are "live" refs).
$ws.Name = "02-10-2026"

# Enter the new production figures for row 16 (cases shipped / pallet count).
$ws.Range("F16").Value = 1836
$ws.Range("M16").Value = 1836

# Enter the manual "Amount Paid" figure for row 79.
$ws.Range("P79").Value = 203796

# Update the view state: scroll the frozen pane over and select the updated
# F77:I77 total cell in the bottom (unfrozen) pane, with the sheet itself
# scrolled down so row 65 is the first visible row.
$ws.Application.ActiveWindow.ScrollRow = 65
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F1").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("F2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("F77:I77").Select()
